$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 24 de Octubre de 2020 a las 08:37"

# --- Reorder Montserrat / Islas Malvinas (rows 216-217 swap countries) ---
$ws.Cells.Item(216, 1).Value = "Montserrat"
$ws.Cells.Item(217, 1).Value = "Islas Malvinas"

# --- Row 26: Ucrania ---
$ws.Cells.Item(26, 2).Value = 337410
$ws.Cells.Item(26, 3).Value = 7014
$ws.Cells.Item(26, 4).Value = 139755
$ws.Cells.Item(26, 5).Value = 191366
$ws.Cells.Item(26, 7).Value = 125
$ws.Cells.Item(26, 8).Value = 6289

# --- Row 63: Uzbekistan ---
$ws.Cells.Item(63, 2).Value = 64811
$ws.Cells.Item(63, 3).Value = 87
$ws.Cells.Item(63, 4).Value = 62033
$ws.Cells.Item(63, 5).Value = 2236

# --- Row 85: El Salvador ---
$ws.Cells.Item(85, 5).Value = 3350
$ws.Cells.Item(85, 7).Value = 4
$ws.Cells.Item(85, 8).Value = 944

# --- Row 122: Sri Lanka ---
$ws.Cells.Item(122, 5).Value = 3494
$ws.Cells.Item(122, 7).Value = 1
$ws.Cells.Item(122, 8).Value = 15

# --- Row 178: Taiwan ---
$ws.Cells.Item(178, 2).Value = 550
$ws.Cells.Item(178, 3).Value = 2
$ws.Cells.Item(178, 5).Value = 46

# --- Row 216: (now Montserrat) activos/muertes swap ---
$ws.Cells.Item(216, 4).Value = 12
$ws.Cells.Item(216, 8).Value = 1

# --- Row 217: (now Islas Malvinas) activos/muertes swap ---
$ws.Cells.Item(217, 4).Value = 13
$ws.Cells.Item(217, 8).Value = 0
